# Update gh-pages to output generated at 456a3b4
# F-column ("想去人数" / interest count) values changed for a handful of
# rows. The same events appear both on the "展览" sheet (their own sheet)
# and again on the aggregated "全部类型" sheet, so both copies get updated.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 302
$ws1.Range("F12").Value = 13298
$ws1.Range("F16").Value = 5473
$ws1.Range("F18").Value = 36

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F21").Value = 302
$ws4.Range("F34").Value = 13298
$ws4.Range("F39").Value = 5473
$ws4.Range("F41").Value = 36
